# Generate Report for Handoff
# Flips the localization status of the single tracked file from
# "Handed back: in sync with en-US" to "Ready for handoff" on every
# sheet, refreshes the handoff-generation timestamps to match, and
# shrinks the now-narrower "Status" columns to fit the new text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"
$newWidth  = 17.2159881591797

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-20 15:01:27"
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-20 15:01:24"
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet ----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-20 15:01:27"
$dede.Columns.Item(3).ColumnWidth = $newWidth
